$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# The shared-string table is rebuilt at save time: any shared string that is
# still referenced by at least one cell keeps its original relative position,
# and any brand-new string value is appended to the table in the order it is
# first assigned to a cell (script order, not sheet order). To land the new
# strings at the exact indices the target workbook expects we must introduce
# them in this precise order:
#   1) "./config_p3388/extreme_rfi.xlsx"  (replaces the orphaned .xlsm path)
#   2) "007"
#   3) "008"
#   4) "009"
#   5) "010"
#   6) "011"
#   7) "012"
# --------------------------------------------------------------------------

# 1) Introduce the corrected extreme-rfi config path first.
$ws.Cells.Item(11, 4).Value = "./config_p3388/extreme_rfi.xlsx"

# 2-7) Introduce the new row IDs, in order, which appends "007".."012".
$ws.Cells.Item(8, 1).Value = "007"
$ws.Cells.Item(9, 1).Value = "008"
$ws.Cells.Item(10, 1).Value = "009"
$ws.Cells.Item(11, 1).Value = "010"
$ws.Cells.Item(12, 1).Value = "011"
$ws.Cells.Item(13, 1).Value = "012"

# --------------------------------------------------------------------------
# Now fill in the rest of the table (all remaining cells reference strings
# that already exist in the workbook, so their write order doesn't matter).
# --------------------------------------------------------------------------

# Row 2: 001 / general_rfi_light (Make flips from 0 to 1)
$ws.Cells.Item(2, 2).Value = 1

# Row 3: 002 -> now general_rfi_light
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = "general_rfi_light"
$ws.Cells.Item(3, 4).Value = "./config_p3388/light_rfi.xlsx"

# Row 4: 003 -> now general_rfi_light
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = "general_rfi_light"
$ws.Cells.Item(4, 4).Value = "./config_p3388/light_rfi.xlsx"

# Row 5: 004 -> now general_rfi_moderate
$ws.Cells.Item(5, 2).Value = 1
$ws.Cells.Item(5, 3).Value = "general_rfi_moderate"
$ws.Cells.Item(5, 4).Value = "./config_p3388/moderate_rfi.xlsx"

# Row 6: 005 -> now general_rfi_moderate
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = "general_rfi_moderate"
$ws.Cells.Item(6, 4).Value = "./config_p3388/moderate_rfi.xlsx"

# Row 7: 006 -> now general_rfi_moderate
$ws.Cells.Item(7, 2).Value = 1
$ws.Cells.Item(7, 3).Value = "general_rfi_moderate"
$ws.Cells.Item(7, 4).Value = "./config_p3388/moderate_rfi.xlsx"

# Row 8: 007 -> general_rfi_severe (new row)
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = "general_rfi_severe"
$ws.Cells.Item(8, 4).Value = "./config_p3388/severe_rfi.xlsx"

# Row 9: 008 -> general_rfi_severe (new row)
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = "general_rfi_severe"
$ws.Cells.Item(9, 4).Value = "./config_p3388/severe_rfi.xlsx"

# Row 10: 009 -> general_rfi_severe (new row)
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = "general_rfi_severe"
$ws.Cells.Item(10, 4).Value = "./config_p3388/severe_rfi.xlsx"

# Row 11: 010 -> general_rfi_extreme (new row, D11 already set above)
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(11, 3).Value = "general_rfi_extreme"

# Row 12: 011 -> general_rfi_extreme (new row)
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(12, 3).Value = "general_rfi_extreme"
$ws.Cells.Item(12, 4).Value = "./config_p3388/extreme_rfi.xlsx"

# Row 13: 012 -> general_rfi_extreme (new row)
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(13, 3).Value = "general_rfi_extreme"
$ws.Cells.Item(13, 4).Value = "./config_p3388/extreme_rfi.xlsx"

# Update the selected cell shown in the saved worksheet view.
$ws.Range("D20").Select()
